$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row just below the header (row 2). This shifts every
# existing data row down by one: the previous last data row (72,
# 2025-11-21) becomes row 73, the previous row 71 (2025-11-22) becomes
# row 72, etc. — exactly the "one more day of history" shift in the diff.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the latest date. Force the cell
# to stay plain text (matching every other date cell in column A) instead
# of letting Excel auto-convert the date-shaped string into a real date
# serial number; then drop back to the Normal style so no stray
# number-format override is left behind on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-01-31"
$ws.Range("A2").Style = "Normal"

# The price columns are constant across the whole history in this sheet,
# so the new row reuses the same figures as every other row.
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
